# Aula15/tabelinha.xlsx - "Ainda alguns errinhos na simulaçao"
# Adds a small "OPERAÇÃO / CODIGO" legend box (H11:I14) explaining the opcode
# simulation, plus a stray "ANTIGO" label at J9, and nudges the view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell text, entered in the same order the original author typed it so the
#    shared-string table comes out in the same order.
# ---------------------------------------------------------------------------
$ws.Range("H11").Value = "OPERAÇÃO"
$ws.Range("H12").Value = "SOMA"

# I12/I13 hold the numeric-looking codes "00"/"01" as TEXT (quote-prefixed),
# so format them (center + middle, like the rest of the table) before typing
# the apostrophe-prefixed value - this is also what creates the first new
# cell style.
$fmtSrc1 = $ws.Range("E3")
$fmtSrc1.Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = "'00"
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = "'01"

$ws.Range("H14").Value = "USE O FUNCT"

# I11 ("CODIGO") is bold, like H11 will become - format it now (second new
# cell style), then fill in the rest of the plain text.
$fmtSrc2 = $ws.Range("E3")
$fmtSrc2.Copy()
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I11").Font.Bold = $true
$ws.Range("I11").Value = "CODIGO"

$ws.Range("J9").Value = "ANTIGO"

# ---------------------------------------------------------------------------
# 2) Formatting for the remaining cells in the legend box.
# ---------------------------------------------------------------------------

# H11 matches I11's bold/centered style.
$fmtSrc3 = $ws.Range("E3")
$fmtSrc3.Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Font.Bold = $true

# H12/H13/H14 are centered (horizontally only, not bold) - third new style.
$fmtSrc4 = $ws.Range("D3")
$fmtSrc4.Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Font.Bold = $false

$ws.Range("H13").Value = "SUB"
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Font.Bold = $false

$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Font.Bold = $false

# J9 and I14 reuse the plain center/middle style used throughout the table.
$fmtSrc5 = $ws.Range("E3")
$fmtSrc5.Copy()
$ws.Range("J9").PasteSpecial(-4122)

$ws.Range("I14").Value = 10
$ws.Range("I14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Column F is a touch narrower now that content settled down.
# ---------------------------------------------------------------------------
$ws.Columns("F").AutoFit()

# ---------------------------------------------------------------------------
# 4) Selection moved to N4 as the last thing the author did before saving.
# ---------------------------------------------------------------------------
$ws.Range("N4").Select()
